$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $ok = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        throw "Replace-Text: anchor not found: $old"
    }
}

function Append-After($anchor, $textToAppend) {
    # Find the anchor text, collapse the found range to its end, then
    # insert the new text right after it (it naturally inherits the
    # surrounding run formatting).
    $r = $d.Content
    $ok = $r.Find.Execute($anchor, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Append-After: anchor not found: $anchor"
    }
    $r.Collapse(0)
    $r.InsertAfter($textToAppend)
}

# ------------------------------------------------------------------
# Title / byline / contact block
# ------------------------------------------------------------------
Replace-Text "Celestial Symphony: The Rhythms of the Universe" "The Allure of History: A Journey Through Time"
Replace-Text "Isabella Maxwell" "Alex Watson"
Replace-Text "isabella" "alex"
Replace-Text "maxwell@astronomycenter" "watson@eduworld"

# ------------------------------------------------------------------
# First body paragraph
# ------------------------------------------------------------------
Replace-Text "Within the vast canvas of the cosmic tapestry, celestial bodies engage in an intricate dance, governed by the harmonious laws of gravitation" "History beckons us, like an alluring whisper from times gone by"
Replace-Text "From the gentle waltz of our solar system's planets to the whirling dervishes of distant galaxies, the cosmos pulsates with an unseen symphony" "It is a kaleidoscope of human experiences, triumphs and follies, wisdom and folly, painted on the canvas of centuries"
Replace-Text "This symphony is a symphony of motion, dictated by the gravitational forces that bind celestial bodies together" "With each chapter, it holds a mirror to our present, shedding light on our origins, shaping our identities, and guiding us into the future"
Replace-Text "It is a symphony of time, measured by the ebb and flow of stellar ages, the rise and fall of civilizations, and the birth and death of stars" "In its vast expanse, history unveils a tapestry woven with countless threads, each representing the lives of individuals who have shaped our world"
Replace-Text "Every celestial object, from the smallest comet to the grandest supermassive black hole, plays a role in this cosmic orchestra, contributing to the intricate composition that orchestrates the universe" "Like intricate brushstrokes, their actions, decisions, and sacrifices add color and texture, creating a vibrant masterpiece that reveals the human capacity for both great achievements and devastating failures"

# New sentence inserted after the "failures" sentence (and its trailing period)
Append-After "Like intricate brushstrokes, their actions, decisions, and sacrifices add color and texture, creating a vibrant masterpiece that reveals the human capacity for both great achievements and devastating failures" "."
Append-After "Like intricate brushstrokes, their actions, decisions, and sacrifices add color and texture, creating a vibrant masterpiece that reveals the human capacity for both great achievements and devastating failures." " Amidst the ebb and flow of civilizations and empires, we discover the timeless struggles of humanity - the quest for power, the pursuit of justice, the yearning for freedom"

Replace-Text "The study of these celestial motions, known as celestial mechanics, delves into the underlying principles that govern the dynamics of the universe" "Furthermore, history teaches us the art of empathy and perspective"
Replace-Text "Scientists, like maestros of the universe, analyze the ballet of planets, the pirouette of stars, and the majestic procession of galaxies" "As we journey through the annals of time, we encounter diverse cultures, beliefs, and ways of life"
Replace-Text "Through this meticulous examination, they unravel the mysteries of the cosmos, revealing its hidden harmonies and unlocking its secrets" "We learn to appreciate the richness of human existence, and we begin to understand why people think, feel, and act as they do"

# New sentence inserted after the "...act as they do" sentence (and its trailing period)
Append-After "We learn to appreciate the richness of human existence, and we begin to understand why people think, feel, and act as they do" "."
Append-After "We learn to appreciate the richness of human existence, and we begin to understand why people think, feel, and act as they do." " This understanding fosters tolerance, compassion, and the realization that we are all part of a shared human story"

Write-Output "stage1+2 ok: $($d.Content.Text.Substring(0, 200))"
